# Update the dSF (column F) values for a handful of rows after a data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -1
    5  = -6
    8  = 0
    9  = 0
    11 = -2
    14 = 3
    17 = 1
    25 = 0
    27 = -1
    28 = -4
    31 = -2
    32 = 4
    34 = -2
    35 = 5
    36 = -2
    37 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
